$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 7857
$ws.Range("I20").Value = 6339.125
$ws.Range("J20").Value = 20000
$ws.Range("K20").Value = 6339.125
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = -6109.125
$ws.Range("N20").Value = -20460
$ws.Range("H29").Value = 3804.8
$ws.Range("I29").Value = 95.5
$ws.Range("J29").Value = 18642
$ws.Range("K29").Value = 286.5
$ws.Range("L29").Value = 55926
$ws.Range("M29").Value = -5.5
$ws.Range("N29").Value = -56488
$ws.Range("H35").Value = 7857
$ws.Range("I35").Value = 6339.125
$ws.Range("J35").Value = 20000
$ws.Range("K35").Value = 6339.125
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = -5960.125
$ws.Range("N35").Value = -20758
$ws.Range("H39").Value = 270.5263
$ws.Range("I39").Value = 17.533333
$ws.Range("J39").Value = 1219.25
$ws.Range("K39").Value = 52.599999
$ws.Range("L39").Value = 3657.75
$ws.Range("M39").Value = 243.400001
$ws.Range("N39").Value = -4249.75
$ws.Range("H53").Value = 155
$ws.Range("I53").Value = 126.5
$ws.Range("K53").Value = 126.5
$ws.Range("M53").Value = 510.5
$ws.Range("H62").Value = 3345
$ws.Range("J62").Value = 4690
$ws.Range("L62").Value = 4690
$ws.Range("N62").Value = -5938
$ws.Range("H64").Value = 15004250
$ws.Range("J64").Value = 5124.75
$ws.Range("L64").Value = 5124.75
$ws.Range("N64").Value = -5620.75
$ws.Range("H65").Value = 3345
$ws.Range("J65").Value = 4690
$ws.Range("L65").Value = 23450
$ws.Range("N65").Value = -29690
$ws.Range("H67").Value = 15004250
$ws.Range("J67").Value = 5124.75
$ws.Range("L67").Value = 5124.75
$ws.Range("N67").Value = -6840.75
$ws.Range("H74").Value = 6375
$ws.Range("I74").Value = 6000
$ws.Range("K74").Value = 6000
$ws.Range("M74").Value = -5064
$ws.Range("H77").Value = 6375
$ws.Range("I77").Value = 6000
$ws.Range("K77").Value = 30000
$ws.Range("M77").Value = -25320
$ws.Range("H100").Value = 1658.7
$ws.Range("I100").Value = 850.75
$ws.Range("K100").Value = 850.75
$ws.Range("M100").Value = -309.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1032.9231
$ws.Range("I32").Value = 759.9
$ws.Range("J32").Value = 1943
$ws.Range("K32").Value = 759.9
$ws.Range("L32").Value = 1943
$ws.Range("M32").Value = -472.9
$ws.Range("N32").Value = -2517
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 10000
$ws.Range("K74").Value = 10000
$ws.Range("M74").Value = -9126
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 10000
$ws.Range("K77").Value = 50000
$ws.Range("M77").Value = -45632
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 207
$ws.Range("I22").Value = 207
$ws.Range("K22").Value = 207
$ws.Range("M22").Value = -34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 330.6154
$ws.Range("J7").Value = 257.25
$ws.Range("L7").Value = 257.25
$ws.Range("N7").Value = -483.25
$ws.Range("H16").Value = 727
$ws.Range("I16").Value = 727
$ws.Range("K16").Value = 727
$ws.Range("M16").Value = -440
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H97").Value = 70000
$ws.Range("J97").Value = 70000
$ws.Range("L97").Value = 70000
$ws.Range("N97").Value = -71982
$ws.Range("H107").Value = 573.5
$ws.Range("I107").Value = 288.2
$ws.Range("K107").Value = 288.2
$ws.Range("M107").Value = 1631.8
$ws.Range("H113").Value = 727
$ws.Range("I113").Value = 727
$ws.Range("K113").Value = 727
$ws.Range("M113").Value = 1443
$ws.Range("H132").Value = 3333.5334
$ws.Range("I132").Value = 2845.3076
$ws.Range("K132").Value = 8535.9228
$ws.Range("M132").Value = -6005.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1100
$ws.Range("I14").Value = 1100
$ws.Range("K14").Value = 3300
$ws.Range("M14").Value = -3127
$ws.Range("H34").Value = 1339.8
$ws.Range("I34").Value = 324.25
$ws.Range("J34").Value = 2016.8334
$ws.Range("K34").Value = 972.75
$ws.Range("L34").Value = 6050.5002
$ws.Range("M34").Value = -888.75
$ws.Range("N34").Value = -6218.5002
$ws.Range("H55").Value = 2611.6667
$ws.Range("J55").Value = 4835
$ws.Range("L55").Value = 14505
$ws.Range("N55").Value = -14859
$ws.Range("H70").Value = 1312
$ws.Range("I70").Value = 1312
$ws.Range("K70").Value = 3936
$ws.Range("M70").Value = -3621
$ws.Range("H73").Value = 1312
$ws.Range("I73").Value = 1312
$ws.Range("K73").Value = 3936
$ws.Range("M73").Value = -2844
$ws.Range("H92").Value = 2225
$ws.Range("I92").Value = 1633.3334
$ws.Range("K92").Value = 4900.0002
$ws.Range("M92").Value = -3652.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5960
$ws.Range("I132").Value = 5325
$ws.Range("K132").Value = 15975
$ws.Range("M132").Value = -13445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8484.583000000001
$ws.Range("I20").Value = 5452.5
$ws.Range("J20").Value = 11516.667
$ws.Range("K20").Value = 5452.5
$ws.Range("L20").Value = 11516.667
$ws.Range("M20").Value = -5226.5
$ws.Range("N20").Value = -11968.667
$ws.Range("H21").Value = 5085.3335
$ws.Range("I21").Value = 5006
$ws.Range("J21").Value = 5125
$ws.Range("K21").Value = 5006
$ws.Range("L21").Value = 5125
$ws.Range("M21").Value = -4832
$ws.Range("N21").Value = -5473
$ws.Range("H24").Value = 15162.2
$ws.Range("I24").Value = 9505.5
$ws.Range("J24").Value = 18933.334
$ws.Range("K24").Value = 9505.5
$ws.Range("L24").Value = 18933.334
$ws.Range("M24").Value = -9162.5
$ws.Range("N24").Value = -19619.334
$ws.Range("H33").Value = 514400
$ws.Range("H93").Value = 21600
$ws.Range("I93").Value = 21600
$ws.Range("K93").Value = 21600
$ws.Range("M93").Value = -20352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6137.7
$ws.Range("I126").Value = 5486.3335
$ws.Range("K126").Value = 16459.0005
$ws.Range("M126").Value = -13989.0005
